$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 32  # was 31
$ws.Range("F8").Value = 64  # was 63
$ws.Range("F10").Value = 10030  # was 10027
$ws.Range("F15").Value = 1940  # was 1939
$ws.Range("F27").Value = 49  # was 48
$ws.Range("F30").Value = 2677  # was 2675
$ws.Range("F31").Value = 923  # was 922
$ws.Range("F32").Value = 610  # was 611
$ws.Range("F36").Value = 473  # was 471
$ws.Range("F37").Value = 200  # was 199
$ws.Range("F39").Value = 1185  # was 1184
$ws.Range("F40").Value = 197  # was 196
$ws.Range("F42").Value = 56  # was 55
$ws.Range("F43").Value = 100  # was 101
$ws.Range("F44").Value = 113  # was 110
$ws.Range("F46").Value = 4033  # was 4034
$ws.Range("F47").Value = 49  # was 48

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 4035  # was 4034

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 32  # was 31
$ws.Range("F9").Value = 64  # was 63
$ws.Range("F11").Value = 10030  # was 10027
$ws.Range("F16").Value = 1940  # was 1939
$ws.Range("F26").Value = 4035  # was 4034
$ws.Range("F29").Value = 49  # was 48
$ws.Range("F32").Value = 2677  # was 2675
$ws.Range("F33").Value = 923  # was 922
$ws.Range("F36").Value = 610  # was 611
$ws.Range("F38").Value = 473  # was 471
$ws.Range("F39").Value = 200  # was 199
$ws.Range("F40").Value = 56  # was 55
$ws.Range("F41").Value = 100  # was 101
$ws.Range("F42").Value = 113  # was 110
$ws.Range("F44").Value = 4033  # was 4034
$ws.Range("F48").Value = 49  # was 48
